# Applies the "Thu Oct 19 17:54:39 UTC 2023" cryptos-list refresh:
# per-row Price (col D) / Volume(1h) (col E) updates, plus a full
# row-51 coin swap (BabyDogeCoin -> Cronos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2: '28.523.11' -> '28.529.53'
$ws.Range("D2").Value = '28.529.53'

# Row 3: D3: '1.561.92' -> '1.560.82'; E3: '  -0.59%  ' -> '  -0.64%  '
$ws.Range("D3").Value = '1.560.82'
$ws.Range("E3").Value = '  -0.64%  '

# Row 4: D4: '1.00' -> '0.999'; E4: '  -0.12%  ' -> '  -0.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.37%  '

# Row 5: D5: '210.31' -> '210.38'; E5: '  -0.73%  ' -> '  -0.72%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.38'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.72%  '

# Row 6: E6: '  -1.33%  ' -> '  -1.26%  '
$ws.Range("E6").Value = '  -1.26%  '

# Row 7: E7: '  -0.07%  ' -> '  -0.20%  '
$ws.Range("E7").Value = '  -0.20%  '

# Row 8: D8: '24.76' -> '24.84'; E8: '  +4.56%  ' -> '  +4.96%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.84'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.96%  '

# Row 9: E9: '  -1.12%  ' -> '  -1.16%  '
$ws.Range("E9").Value = '  -1.16%  '

# Row 10: E10: '  -0.53%  ' -> '  -0.48%  '
$ws.Range("E10").Value = '  -0.48%  '

# Row 11: D11: '0.0896' -> '0.0895'; E11: '  +0.31%  ' -> '  +0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0895'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.14%  '

# Row 12: D12: '1.787.72' -> '1.785.82'; E12: '  -0.48%  ' -> '  -0.54%  '
$ws.Range("D12").Value = '1.785.82'
$ws.Range("E12").Value = '  -0.54%  '

# Row 13: D13: '1.560.86' -> '1.587.03'; E13: '  -0.69%  ' -> '  +0.96%  '
$ws.Range("D13").Value = '1.587.03'
$ws.Range("E13").Value = '  +0.96%  '

# Row 14: D14: '28.565.40' -> '28.554.46'; E14: '  +0.59%  ' -> '  +0.52%  '
$ws.Range("D14").Value = '28.554.46'
$ws.Range("E14").Value = '  +0.52%  '

# Row 15: E15: '  -0.54%  ' -> '  -0.40%  '
$ws.Range("E15").Value = '  -0.40%  '

# Row 16: E16: '  -1.90%  ' -> '  -1.71%  '
$ws.Range("E16").Value = '  -1.71%  '

# Row 17: D17: '61.12' -> '61.13'; E17: '  -0.81%  ' -> '  -0.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.13'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.78%  '

# Row 18: D18: '229.36' -> '229.72'; E18: '  +0.01%  ' -> '  +0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.72'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.40%  '

# Row 19: D19: '7.35' -> '7.34'; E19: '  -0.49%  ' -> '  -0.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.34'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.64%  '

# Row 20: D20: '0.0₃0673' -> '0.0₃0675'; E20: '  -1.49%  ' -> '  -1.37%  '
$ws.Range("D20").Value = '0.0₃0675'
$ws.Range("E20").Value = '  -1.37%  '

# Row 21: E21: '  -0.11%  ' -> '  -0.21%  '
$ws.Range("E21").Value = '  -0.21%  '

# Row 23: D23: '8.95' -> '8.96'; E23: '  -0.74%  ' -> '  -0.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.96'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.43%  '

# Row 24: D24: '2.07' -> '2.08'; E24: '  +1.01%  ' -> '  +1.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.08'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.16%  '

# Row 25: D25: '150.17' -> '150.15'; E25: '  -0.80%  ' -> '  -0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.15'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.70%  '

# Row 26: D26: '14.76' -> '14.77'; E26: '  -0.98%  ' -> '  -0.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.77'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.87%  '

# Row 27: E27: '  -0.32%  ' -> '  -0.27%  '
$ws.Range("E27").Value = '  -0.27%  '

# Row 28: E28: '  -0.05%  ' -> '  -0.18%  '
$ws.Range("E28").Value = '  -0.18%  '

# Row 29: E29: '  -2.52%  ' -> '  -2.40%  '
$ws.Range("E29").Value = '  -2.40%  '

# Row 30: E30: '  -4.55%  ' -> '  -4.50%  '
$ws.Range("E30").Value = '  -4.50%  '

# Row 31: E31: '  -1.69%  ' -> '  -2.00%  '
$ws.Range("E31").Value = '  -2.00%  '

# Row 32: E32: '  -0.81%  ' -> '  -0.67%  '
$ws.Range("E32").Value = '  -0.67%  '

# Row 33: D33: '1.386.70' -> '1.387.51'; E33: '  +0.45%  ' -> '  +0.65%  '
$ws.Range("D33").Value = '1.387.51'
$ws.Range("E33").Value = '  +0.65%  '

# Row 34: D34: '2.97' -> '2.96'; E34: '  -4.24%  ' -> '  -4.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.96'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.31%  '

# Row 35: E35: '  -2.99%  ' -> '  -2.86%  '
$ws.Range("E35").Value = '  -2.86%  '

# Row 36: E36: '  -2.10%  ' -> '  -1.84%  '
$ws.Range("E36").Value = '  -1.84%  '

# Row 37: D37: '2.69' -> '2.68'; E37: '  +1.52%  ' -> '  +1.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.68'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.23%  '

# Row 38: E38: '  -2.63%  ' -> '  -2.69%  '
$ws.Range("E38").Value = '  -2.69%  '

# Row 39: E39: '  -1.31%  ' -> '  -1.25%  '
$ws.Range("E39").Value = '  -1.25%  '

# Row 40: E40: '  +2.12%  ' -> '  +2.18%  '
$ws.Range("E40").Value = '  +2.18%  '

# Row 41: E41: '  -0.79%  ' -> '  -1.00%  '
$ws.Range("E41").Value = '  -1.00%  '

# Row 42: E42: '  -0.12%  ' -> '  -0.19%  '
$ws.Range("E42").Value = '  -0.19%  '

# Row 43: D43: '0.770' -> '0.771'; E43: '  -2.01%  ' -> '  -2.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.771'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.04%  '

# Row 44: E44: '  -2.85%  ' -> '  -2.86%  '
$ws.Range("E44").Value = '  -2.86%  '

# Row 45: D45: '63.70' -> '63.71'; E45: '  +2.25%  ' -> '  +2.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.71'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.38%  '

# Row 46: E46: '  -2.48%  ' -> '  -2.37%  '
$ws.Range("E46").Value = '  -2.37%  '

# Row 47: D47: '1.699.77' -> '1.697.15'; E47: '  -0.54%  ' -> '  -0.65%  '
$ws.Range("D47").Value = '1.697.15'
$ws.Range("E47").Value = '  -0.65%  '

# Row 48: E48: '  -5.40%  ' -> '  -5.39%  '
$ws.Range("E48").Value = '  -5.39%  '

# Row 49: D49: '85.04' -> '85.06'; E49: '  -0.32%  ' -> '  -0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.06'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.30%  '

# Row 50: D50: '42.84' -> '42.98'; E50: '  +5.99%  ' -> '  +6.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '42.98'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.33%  '

# Row 51: B51: 'BabyDogeCoin' -> 'Cronos'; C51: 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge' -> 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; D51: '0.0₆0100' -> '0.0512'; E51: '  +0.69%  ' -> '  -0.58%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0512'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.58%  '
